$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Edit "Sheet1" (the ToDo checklist) in place: blank out the answer cells
#    (keep their formatting), drop the two duplicate "aaa" rows, and turn the
#    top of the sheet into a small "Missing Contents" callout with three new
#    notes.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2-4: move the label text from column B into column A with new wording;
# the old B cells vanish completely (their style matched the column default
# anyway), while column C formatting on row 4 is kept but its value blanked.
$ws1.Range("B2").Clear()
$ws1.Range("A2").Value = "API protocol for all pages"

$ws1.Range("B3").Clear()
$ws1.Range("A3").Value = "Responsible(Flex)"

$ws1.Range("B4").Clear()
$ws1.Range("A4").Value = "Table height( when count is less than one page capacity)"
$ws1.Range("C4").ClearContents()

# Rows 5-10: column C values cleared, formatting preserved.
$ws1.Range("C5").ClearContents()
$ws1.Range("C6").ClearContents()
$ws1.Range("C7").ClearContents()
$ws1.Range("C8").ClearContents()
$ws1.Range("C9").ClearContents()
$ws1.Range("C10").ClearContents()

# Row 11: the B cell disappears entirely (no formatting left behind either);
# C keeps its style but loses its value.
$ws1.Range("B11").Clear()
$ws1.Range("C11").ClearContents()

# Rows 12-13: both columns cleared, formatting preserved.
$ws1.Range("B12").ClearContents()
$ws1.Range("C12").ClearContents()
$ws1.Range("B13").ClearContents()
$ws1.Range("C13").ClearContents()

# Rows 14-21: column C values cleared.
$ws1.Range("C14").ClearContents()
$ws1.Range("C15").ClearContents()
$ws1.Range("C16").ClearContents()
$ws1.Range("C17").ClearContents()
$ws1.Range("C18").ClearContents()
$ws1.Range("C19").ClearContents()
$ws1.Range("C20").ClearContents()
$ws1.Range("C21").ClearContents()

# Row 22 ("aaa") is removed completely.
$ws1.Range("B22").Clear()

# Rows 25-27: clear values, keep formatting.
$ws1.Range("B25").ClearContents()
$ws1.Range("C25").ClearContents()
$ws1.Range("C26").ClearContents()
$ws1.Range("B27").ClearContents()

# Row 28 ("aaa" duplicate) is removed completely.
$ws1.Range("B28").Clear()

# Row 29: clear value, keep formatting.
$ws1.Range("B29").ClearContents()

# Row 32: B disappears entirely, C keeps formatting only.
$ws1.Range("B32").Clear()
$ws1.Range("C32").ClearContents()

# Rows 33-35: clear values, keep formatting.
$ws1.Range("C33").ClearContents()
$ws1.Range("C34").ClearContents()
$ws1.Range("B35").ClearContents()

# New header row at the very top: a yellow "Missing Contents" callout.
$ws1.Range("A1").Value = "Missing Contents"
$ws1.Range("A1").Interior.Color = 65535

# Column A becomes the wide "label" column now that it holds real text.
$ws1.Columns("A").ColumnWidth = 45.5

# ---------------------------------------------------------------------------
# 2) Re-order the tabs: "Sheet2" (the API/component tracking table) moves in
#    front of "Sheet1", and "Sheet1" becomes the active tab.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Move($wb.Worksheets.Item(1))

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Select()
